# "Store the 4 stations into Mongo Atlas (question 1)"
# The data-catalog sheet's field-name column (A2:A8) is re-cased to a
# Title-Case / acronym style ("name" -> "Name", "tpe" -> "TPE", ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Name"
$ws.Range("A3").Value = "Coordinates"
$ws.Range("A4").Value = "Size"
$ws.Range("A5").Value = "TPE"
$ws.Range("A6").Value = "Available"
$ws.Range("A7").Value = "City"
$ws.Range("A8").Value = "Municipality"

# The author resized the columns (A now has its own explicit width, B/C
# share a width, D got wider, E stayed about the same) and left the
# cursor parked on B15 before saving.
$ws.Columns.Item(1).ColumnWidth = 14.833333333333332
$ws.Range("B1:C1").EntireColumn.ColumnWidth = 27
$ws.Columns.Item(4).ColumnWidth = 33.5
$ws.Columns.Item(5).ColumnWidth = 16.833333333333336

$ws.Range("B15").Select()
